$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Range("N14").Value = -90

# Row 15
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -33.333333333333
$ws.Range("N15").Value = -67.1875

# Row 16
$ws.Range("C16").Value = 3
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 120
$ws.Range("I16").Value = 110
$ws.Range("K16").Value = 8.910891089108
$ws.Range("L16").Value = -4.347826086956
$ws.Range("M16").Value = -58.174904942965
$ws.Range("N16").Value = -87.442922374429

# Row 17
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 23
$ws.Range("H17").Value = -4.166666666666
$ws.Range("I17").Value = 276
$ws.Range("J17").Value = 312
$ws.Range("K17").Value = -11.538461538461
$ws.Range("L17").Value = -11.821086261980
$ws.Range("M17").Value = 2.985074626865
$ws.Range("N17").Value = -51.063829787234

# Row 18
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -71.428571428571
$ws.Range("I18").Value = 61
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = -39
$ws.Range("L18").Value = -33.695652173913
$ws.Range("M18").Value = -79.037800687285
$ws.Range("N18").Value = -92.738095238095

# Row 19
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 30
$ws.Range("H19").Value = -9.090909090909
$ws.Range("I19").Value = 317
$ws.Range("J19").Value = 341
$ws.Range("K19").Value = -7.038123167155
$ws.Range("L19").Value = 19.622641509434
$ws.Range("M19").Value = -33.820459290187
$ws.Range("N19").Value = -90.287990196078

# Row 20
$ws.Range("C20").Value = 7
$ws.Range("E20").Value = 600
$ws.Range("F20").Value = 22
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 166
$ws.Range("J20").Value = 183
$ws.Range("K20").Value = -9.289617486338
$ws.Range("L20").Value = 20.289855072463
$ws.Range("M20").Value = -22.790697674418
$ws.Range("N20").Value = -87.376425855513

# Row 21
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 47.058823529411
$ws.Range("F21").Value = 92
$ws.Range("G21").Value = 90
$ws.Range("H21").Value = 2.222222222222
$ws.Range("I21").Value = 954
$ws.Range("J21").Value = 1078
$ws.Range("K21").Value = -11.502782931354
$ws.Range("L21").Value = 0.104931794333
$ws.Range("M21").Value = -38.451612903225
$ws.Range("N21").Value = -86.279303897598

# Row 23
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F23").PasteSpecial(-4122)
$ws.Range("H23").Value = -100
$ws.Range("M23").Value = 0

# Row 24
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -33.333333333333
$ws.Range("F24").Value = 66
$ws.Range("G24").Value = 114
$ws.Range("H24").Value = -42.105263157894
$ws.Range("I24").Value = 902
$ws.Range("J24").Value = 1021
$ws.Range("K24").Value = -11.655239960822
$ws.Range("L24").Value = 24.242424242424
$ws.Range("M24").Value = 12.609238451935

# Row 25
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -14.285714285714
$ws.Range("F25").Value = 54
$ws.Range("G25").Value = 38
$ws.Range("H25").Value = 42.105263157894
$ws.Range("I25").Value = 515
$ws.Range("J25").Value = 407
$ws.Range("K25").Value = 26.535626535626
$ws.Range("L25").Value = 38.069705093833
$ws.Range("M25").Value = -18.897637795275

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("F26").Value = 6
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 100
$ws.Range("I26").Value = 37
$ws.Range("K26").Value = -15.909090909090
$ws.Range("L26").Value = 0

# Row 27
$ws.Range("C27").Value = 2
$ws.Range("I14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -33.333333333333
$ws.Range("F27").Value = 2
$ws.Range("I14").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -66.666666666666
$ws.Range("I27").Value = 38
$ws.Range("J27").Value = 43
$ws.Range("K27").Value = -11.627906976744
$ws.Range("L27").Value = 0

# Row 28
$ws.Range("D28").Value = 2
$ws.Range("I14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("G28").Value = 3
$ws.Range("J28").Value = 30
$ws.Range("K28").Value = -46.666666666666
$ws.Range("M28").Value = -66.666666666666
$ws.Range("N28").Value = -86.440677966101

# Row 29
$ws.Range("D29").Value = 2
$ws.Range("I14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("G29").Value = 3
$ws.Range("J29").Value = 24
$ws.Range("K29").Value = -54.166666666666
$ws.Range("M29").Value = -71.794871794871
$ws.Range("N29").Value = -89.523809523809

$ws.Application.CutCopyMode = $false

# Shared-string text fixups (report header)
$ws.Range("A8").Value = "Volume 30   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  10/2/2023  Through  10/8/2023"

